$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.029376799356197
$ws.Range("D2").Value = 1.033094477126458
$ws.Range("E2").Value = 1.038765405128893
$ws.Range("F2").Value = 1.048752612011598
$ws.Range("I2").Value = 1.035446552766692
$ws.Range("J2").Value = 1.034523914377592
$ws.Range("K2").Value = 1.035897815932917
$ws.Range("L2").Value = 1.041552506389207
$ws.Range("M2").Value = 1.051511580072573
$ws.Range("N2").Value = 1.015523357451605
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.030166194290693
$ws.Range("D3").Value = 1.033667725293768
$ws.Range("E3").Value = 1.039540673569085
$ws.Range("F3").Value = 1.049777376184659
$ws.Range("I3").Value = 1.035618480492842
$ws.Range("J3").Value = 1.034955134405948
$ws.Range("K3").Value = 1.036280602726234
$ws.Range("L3").Value = 1.042137947312219
$ws.Range("M3").Value = 1.052347902863642
$ws.Range("N3").Value = 1.015666406153588
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.030677412003788
$ws.Range("D4").Value = 1.034038896022996
$ws.Range("E4").Value = 1.040043167501797
$ws.Range("F4").Value = 1.050441853987587
$ws.Range("I4").Value = 1.035728514568565
$ws.Range("J4").Value = 1.035233915036383
$ws.Range("K4").Value = 1.036527830289363
$ws.Range("L4").Value = 1.042516950260741
$ws.Range("M4").Value = 1.052889837088291
$ws.Range("N4").Value = 1.015758869003269
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.030892428374054
$ws.Range("D5").Value = 1.034194992124581
$ws.Range("E5").Value = 1.040254616169652
$ws.Range("F5").Value = 1.050721530826112
$ws.Range("I5").Value = 1.035774481324573
$ws.Range("J5").Value = 1.035351053970645
$ws.Range("K5").Value = 1.036631653171344
$ws.Range("L5").Value = 1.042676325699563
$ws.Range("M5").Value = 1.05311785084507
$ws.Range("N5").Value = 1.015797716212017
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.030928536411195
$ws.Range("D6").Value = 1.034221204576137
$ws.Range("E6").Value = 1.04029013102301
$ws.Range("F6").Value = 1.050768509123212
$ws.Range("I6").Value = 1.035782182224805
$ws.Range("J6").Value = 1.035370718534557
$ws.Range("K6").Value = 1.036649078911724
$ws.Range("L6").Value = 1.042703087996215
$ws.Range("M6").Value = 1.053156146139679
$ws.Range("N6").Value = 1.015804237399598
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.030680284670426
$ws.Range("D7").Value = 1.034040981569169
$ws.Range("E7").Value = 1.040045992104652
$ws.Range("F7").Value = 1.050445589747121
$ws.Range("I7").Value = 1.035729129925248
$ws.Range("J7").Value = 1.035235480491528
$ws.Range("K7").Value = 1.036529218015755
$ws.Range("L7").Value = 1.042519079677446
$ws.Range("M7").Value = 1.05289288309578
$ws.Range("N7").Value = 1.015759388177083
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.029643489794756
$ws.Range("D8").Value = 1.033288158124142
$ws.Range("E8").Value = 1.039027235074687
$ws.Range("F8").Value = 1.049098647712364
$ws.Range("I8").Value = 1.035504907582472
$ws.Range("J8").Value = 1.034669697757948
$ws.Range("K8").Value = 1.036027275272199
$ws.Range("L8").Value = 1.041750320086494
$ws.Range("M8").Value = 1.051794057386108
$ws.Range("N8").Value = 1.015571721724166
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.02781986296577
$ws.Range("D9").Value = 1.031963507832455
$ws.Range("E9").Value = 1.037238589513838
$ws.Range("F9").Value = 1.046735858338226
$ws.Range("I9").Value = 1.035100530690744
$ws.Range("J9").Value = 1.033670872201578
$ws.Range("K9").Value = 1.035139309636468
$ws.Range("L9").Value = 1.040397132140116
$ws.Range("M9").Value = 1.049863807550775
$ws.Range("N9").Value = 1.015240287980712
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.026606453583973
$ws.Range("D10").Value = 1.031081802004838
$ws.Range("E10").Value = 1.036050653530774
$ws.Range("F10").Value = 1.045167953742466
$ws.Range("I10").Value = 1.034824760740472
$ws.Range("J10").Value = 1.033003818552254
$ws.Range("K10").Value = 1.034545065246852
$ws.Range("L10").Value = 1.039496071135833
$ws.Range("M10").Value = 1.048581112945244
$ws.Range("N10").Value = 1.015018858503526
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.02608160983478
$ws.Range("D11").Value = 1.03070036610385
$ws.Range("E11").Value = 1.035537351414018
$ws.Range("F11").Value = 1.044490782951727
$ws.Range("I11").Value = 1.034703891931697
$ws.Range("J11").Value = 1.03271471416881
$ws.Range("K11").Value = 1.034287228117203
$ws.Range("L11").Value = 1.039106170932461
$ws.Range("M11").Value = 1.048026692554024
$ws.Range("N11").Value = 1.014922869988438
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.025886746849201
$ws.Range("D12").Value = 1.030558737948317
$ws.Range("E12").Value = 1.035346852428118
$ws.Range("F12").Value = 1.044239514966196
$ws.Range("I12").Value = 1.034658777498091
$ws.Range("J12").Value = 1.032607289178094
$ws.Range("K12").Value = 1.034191378194855
$ws.Range("L12").Value = 1.038961385813654
$ws.Range("M12").Value = 1.047820907111716
$ws.Range("N12").Value = 1.014887199759918
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.025928541641943
$ws.Range("D13").Value = 1.030589115191164
$ws.Range("E13").Value = 1.035387707645913
$ws.Range("F13").Value = 1.04429340085632
$ws.Range("I13").Value = 1.034668464566995
$ws.Range("J13").Value = 1.032630333937927
$ws.Range("K13").Value = 1.03421194182684
$ws.Range("L13").Value = 1.038992440828254
$ws.Range("M13").Value = 1.047865041931309
$ws.Range("N13").Value = 1.014894851854189
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.026065500608186
$ws.Range("D14").Value = 1.030688657962124
$ws.Range("E14").Value = 1.035521601347332
$ws.Range("F14").Value = 1.044470007678912
$ws.Range("I14").Value = 1.034700167208944
$ws.Range("J14").Value = 1.032705835169263
$ws.Range("K14").Value = 1.034279306711492
$ws.Range("L14").Value = 1.039094202101624
$ws.Range("M14").Value = 1.048009679176529
$ws.Range("N14").Value = 1.014919921796525
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.02614989717143
$ws.Range("D15").Value = 1.03074999677275
$ws.Range("E15").Value = 1.035604119495192
$ws.Range("F15").Value = 1.044578855923012
$ws.Range("I15").Value = 1.034719671338784
$ws.Range("J15").Value = 1.03275234887011
$ws.Range("K15").Value = 1.034320802175887
$ws.Range("L15").Value = 1.03915690607714
$ws.Range("M15").Value = 1.048098815008119
$ws.Range("N15").Value = 1.01493536613162
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.02664129789501
$ws.Range("D16").Value = 1.031107124142803
$ws.Range("E16").Value = 1.036084742644816
$ws.Range("F16").Value = 1.045212932185032
$ws.Range("I16").Value = 1.034832751739221
$ws.Range("J16").Value = 1.033022999959342
$ws.Range("K16").Value = 1.034562166081244
$ws.Range("L16").Value = 1.039521953229169
$ws.Range("M16").Value = 1.048617929089874
$ws.Range("N16").Value = 1.015025226703598
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.026949694486176
$ws.Range("D17").Value = 1.031331235262933
$ws.Range("E17").Value = 1.036386515743364
$ws.Range("F17").Value = 1.045611139080435
$ws.Range("I17").Value = 1.034903294115746
$ws.Range("J17").Value = 1.03319270191814
$ws.Range("K17").Value = 1.034713427469655
$ws.Range("L17").Value = 1.03975100959932
$ws.Range("M17").Value = 1.04894382305232
$ws.Range("N17").Value = 1.015081565223572
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.027129631937224
$ws.Range("D18").Value = 1.031461988999163
$ws.Range("E18").Value = 1.036562639168685
$ws.Range("F18").Value = 1.045843574416357
$ws.Range("I18").Value = 1.034944299505565
$ws.Range("J18").Value = 1.033291660479601
$ws.Range("K18").Value = 1.034801604905607
$ws.Range("L18").Value = 1.039884639798669
$ws.Range("M18").Value = 1.04913400737266
$ws.Range("N18").Value = 1.015114416071872
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.02719099522652
$ws.Range("D19").Value = 1.031506578274851
$ws.Range("E19").Value = 1.036622710341187
$ws.Range("F19").Value = 1.045922857303451
$ws.Range("I19").Value = 1.034958257381084
$ws.Range("J19").Value = 1.033325398393408
$ws.Range("K19").Value = 1.034831662484047
$ws.Range("L19").Value = 1.03993020852679
$ws.Range("M19").Value = 1.049198871528976
$ws.Range("N19").Value = 1.015125615575713
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.02691660075636
$ws.Range("D20").Value = 1.031307186781772
$ws.Range("E20").Value = 1.036354127547022
$ws.Range("F20").Value = 1.045568397861182
$ws.Range("I20").Value = 1.034895740142897
$ws.Range("J20").Value = 1.033174497155853
$ws.Range("K20").Value = 1.034697203797194
$ws.Range("L20").Value = 1.039726431371936
$ws.Range("M20").Value = 1.048908847785538
$ws.Range("N20").Value = 1.015075521712175
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.026025167181753
$ws.Range("D21").Value = 1.030659343589004
$ws.Range("E21").Value = 1.0354821684409
$ws.Range("F21").Value = 1.04441799408719
$ws.Range("I21").Value = 1.034690837593192
$ws.Range("J21").Value = 1.032683602998713
$ws.Range("K21").Value = 1.034259471568136
$ws.Range("L21").Value = 1.039064234795497
$ws.Range("M21").Value = 1.047967082932939
$ws.Range("N21").Value = 1.014912539756707
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.025465193713138
$ws.Range("D22").Value = 1.030252332858093
$ws.Range("E22").Value = 1.034934884436955
$ws.Range("F22").Value = 1.043696214326896
$ws.Range("I22").Value = 1.03456074414709
$ws.Range("J22").Value = 1.032374734742755
$ws.Range("K22").Value = 1.033983802801375
$ws.Range("L22").Value = 1.038648124426717
$ws.Range("M22").Value = 1.047375832246686
$ws.Range("N22").Value = 1.014809975209355
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.02576199771482
$ws.Range("D23").Value = 1.030468066598499
$ws.Range("E23").Value = 1.035224919297721
$ws.Range("F23").Value = 1.04407869835534
$ws.Range("I23").Value = 1.034629828615178
$ws.Range("J23").Value = 1.032538492452938
$ws.Range("K23").Value = 1.034129982203023
$ws.Range("L23").Value = 1.038868689313584
$ws.Range("M23").Value = 1.047689182046889
$ws.Range("N23").Value = 1.014864355133022
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.026931554218809
$ws.Range("D24").Value = 1.031318053151934
$ws.Range("E24").Value = 1.036368762056812
$ws.Range("F24").Value = 1.045587710259508
$ws.Range("I24").Value = 1.034899153893162
$ws.Range("J24").Value = 1.033182723183928
$ws.Range("K24").Value = 1.034704534733466
$ws.Range("L24").Value = 1.039737537136321
$ws.Range("M24").Value = 1.048924651307005
$ws.Range("N24").Value = 1.015078252547198
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.028290908585207
$ws.Range("D25").Value = 1.032305723438715
$ws.Range("E25").Value = 1.037700212489159
$ws.Range("F25").Value = 1.047345418592022
$ws.Range("I25").Value = 1.035206165409235
$ws.Range("J25").Value = 1.033929303877368
$ws.Range("K25").Value = 1.035369274885491
$ws.Range("L25").Value = 1.040746781895387
$ws.Range("M25").Value = 1.050362100550425
$ws.Range("N25").Value = 1.015326056851311
